$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("zh-cn")
$ws.Range("K8").Value = "2016-08-13 06:53:41"
Write-Output $ws.Range("K8").Value
